$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels in row 1 (N1, O1)
$ws.Range("N1").Value = "Case_0 with Openness"
$ws.Range("O1").Value = "Case_0 with Voiceless"

# Update nationality value in B11
$ws.Range("B11").Value = "Indonesia"

# Swap the N and O column values for data rows 2 through 19
for ($r = 2; $r -le 19; $r++) {
    $nCell = $ws.Cells.Item($r, 14)  # Column N
    $oCell = $ws.Cells.Item($r, 15)  # Column O

    $nVal = $nCell.Value2
    $oVal = $oCell.Value2

    $nCell.Value = $oVal
    $oCell.Value = $nVal
}
